$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume(1h) (E) columns keep their original
# text representation (e.g. "235.50", "1.000", "24.890.81") instead of
# being auto-coerced to numbers by Excel when the text looks numeric.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.890.81"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.634.82"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -6.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -5.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4725"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -6.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2560"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -6.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06059"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06950"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.637.83"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.65"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6049"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -7.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.333"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "72.84"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.892.05"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006561"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.11"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.848.64"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.344"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.551"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.221"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "133.15"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.60%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.382"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -7.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "102.95"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.630"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -8.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.762"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07728"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.529"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9987"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04292"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -8.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.584"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9197"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -7.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5809"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.537"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01537"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9984"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8174"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.97"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.773"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -7.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3683"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.707"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1090"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05194"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.029"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.29%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9998"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9957"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.49%  "
